# Scheduled data refresh: update crafting-profit market-price figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) across all eight job sheets with the latest pulled values.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 6334.931
$ws.Range("I38").Value = 6301
$ws.Range("J38").Value = 6399.4
$ws.Range("K38").Value = 18903
$ws.Range("L38").Value = 19198.2
$ws.Range("M38").Value = -18531
$ws.Range("N38").Value = -19942.2

$ws.Range("H106").Value = 9174.842000000001
$ws.Range("I106").Value = 8613.923000000001
$ws.Range("J106").Value = 10390.167
$ws.Range("K106").Value = 8613.923000000001
$ws.Range("L106").Value = 10390.167
$ws.Range("M106").Value = -7982.923000000001
$ws.Range("N106").Value = -11652.167

$ws.Range("H129").Value = 2396.4443
$ws.Range("J129").Value = 2983.4
$ws.Range("L129").Value = 8950.200000000001
$ws.Range("N129").Value = -18950.2

$ws.Range("H137").Value = 3093.3125
$ws.Range("I137").Value = 2829.2
$ws.Range("K137").Value = 8487.599999999999
$ws.Range("M137").Value = -5937.599999999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2391.389
$ws.Range("J2").Value = 3738.6365
$ws.Range("L2").Value = 3738.6365
$ws.Range("N2").Value = -3964.6365

$ws.Range("H32").Value = 2063.057
$ws.Range("I32").Value = 1718.7412
$ws.Range("K32").Value = 1718.7412
$ws.Range("M32").Value = -1431.7412

$ws.Range("H74").Value = 1631.6786
$ws.Range("I74").Value = 1588.6666
$ws.Range("J74").Value = 1889.75
$ws.Range("K74").Value = 1588.6666
$ws.Range("L74").Value = 1889.75
$ws.Range("M74").Value = -714.6666
$ws.Range("N74").Value = -3637.75

$ws.Range("H77").Value = 1631.6786
$ws.Range("I77").Value = 1588.6666
$ws.Range("J77").Value = 1889.75
$ws.Range("K77").Value = 7943.333000000001
$ws.Range("L77").Value = 9448.75
$ws.Range("M77").Value = -3575.333000000001
$ws.Range("N77").Value = -18184.75

$ws.Range("H102").Value = 7271.32
$ws.Range("I102").Value = 5863
$ws.Range("K102").Value = 5863
$ws.Range("M102").Value = -4241

$ws.Range("H116").Value = 2391.389
$ws.Range("J116").Value = 3738.6365
$ws.Range("L116").Value = 3738.6365
$ws.Range("N116").Value = -8326.636500000001

$ws.Range("H132").Value = 2324.25
$ws.Range("I132").Value = 2687.6875
$ws.Range("K132").Value = 8063.0625
$ws.Range("M132").Value = -5533.0625

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2391.389
$ws.Range("J3").Value = 3738.6365
$ws.Range("L3").Value = 3738.6365
$ws.Range("N3").Value = -3966.6365

$ws.Range("H82").Value = 20748
$ws.Range("I82").Value = 5288.6665
$ws.Range("J82").Value = 51666.668
$ws.Range("K82").Value = 5288.6665
$ws.Range("L82").Value = 51666.668
$ws.Range("M82").Value = -4905.6665
$ws.Range("N82").Value = -52432.668

$ws.Range("H85").Value = 20748
$ws.Range("I85").Value = 5288.6665
$ws.Range("J85").Value = 51666.668
$ws.Range("K85").Value = 5288.6665
$ws.Range("L85").Value = 51666.668
$ws.Range("M85").Value = -3962.6665
$ws.Range("N85").Value = -54318.668

$ws.Range("H99").Value = 3560
$ws.Range("I99").Value = 1612.5
$ws.Range("K99").Value = 1612.5
$ws.Range("M99").Value = -114.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2063.5454
$ws.Range("I2").Value = 2700
$ws.Range("J2").Value = 366.33334
$ws.Range("K2").Value = 2700
$ws.Range("L2").Value = 366.33334
$ws.Range("M2").Value = -2587
$ws.Range("N2").Value = -592.33334

$ws.Range("H132").Value = 2132.3845
$ws.Range("I132").Value = 1851.7916
$ws.Range("K132").Value = 5555.3748
$ws.Range("M132").Value = -3025.3748

$ws.Range("H141").Value = 155999.75
$ws.Range("J141").Value = 155999.75
$ws.Range("L141").Value = 155999.75
$ws.Range("N141").Value = -166359.75

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2174.3684
$ws.Range("J113").Value = 2557.6
$ws.Range("L113").Value = 7672.799999999999
$ws.Range("N113").Value = -12012.8

$ws.Range("H121").Value = 765.9091
$ws.Range("J121").Value = 1100
$ws.Range("L121").Value = 3300
$ws.Range("N121").Value = -5920

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 633499.9
$ws.Range("I10").Value = 1672666.4
$ws.Range("K10").Value = 1672666.4
$ws.Range("M10").Value = -1672497.4

$ws.Range("H21").Value = 9000
$ws.Range("J21").Value = 9000
$ws.Range("L21").Value = 9000
$ws.Range("N21").Value = -9346

$ws.Range("H22").Value = 3792.5
$ws.Range("I22").Value = 2600
$ws.Range("K22").Value = 2600
$ws.Range("M22").Value = -2071

$ws.Range("H30").Value = 9000
$ws.Range("J30").Value = 9000
$ws.Range("L30").Value = 9000
$ws.Range("N30").Value = -9210

$ws.Range("H80").Value = 4407.077
$ws.Range("J80").Value = 5361.75
$ws.Range("L80").Value = 5361.75
$ws.Range("N80").Value = -7357.75

$ws.Range("H83").Value = 4407.077
$ws.Range("J83").Value = 5361.75
$ws.Range("L83").Value = 26808.75
$ws.Range("N83").Value = -36792.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 6137.4443
$ws.Range("I68").Value = 4099.5
$ws.Range("K68").Value = 4099.5
$ws.Range("M68").Value = -3350.5

$ws.Range("H71").Value = 6137.4443
$ws.Range("I71").Value = 4099.5
$ws.Range("K71").Value = 20497.5
$ws.Range("M71").Value = -16753.5

$ws.Range("H93").Value = 3627.182
$ws.Range("I93").Value = 1159.8
$ws.Range("K93").Value = 1159.8
$ws.Range("M93").Value = 88.20000000000005

$ws.Range("H100").Value = 4801.2104
$ws.Range("I100").Value = 3529.2727
$ws.Range("K100").Value = 3529.2727
$ws.Range("M100").Value = -2988.2727

$ws.Range("H136").Value = 16195.5
$ws.Range("I136").Value = 1045.8235
$ws.Range("J136").Value = 29750.475
$ws.Range("K136").Value = 3137.4705
$ws.Range("L136").Value = 89251.42499999999
$ws.Range("M136").Value = -587.4704999999999
$ws.Range("N136").Value = -94351.42499999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 25000
$ws.Range("J54").Value = 25000
$ws.Range("L54").Value = 25000
$ws.Range("N54").Value = -26040

$ws.Range("H132").Value = 2038.0385
$ws.Range("I132").Value = 791.2414
$ws.Range("K132").Value = 2373.7242
$ws.Range("M132").Value = 156.2757999999999

$ws.Range("H136").Value = 1741.6538
$ws.Range("I136").Value = 1315.619
$ws.Range("K136").Value = 3946.857
$ws.Range("M136").Value = -1396.857
